$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.253.63"
$ws.Range("E2").Value = "  -1.88%  "
$ws.Range("D3").Value = "1.875.22"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'307.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.5200"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.47%  "
$ws.Range("D8").Value = "'0.3751"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.51%  "
$ws.Range("D9").Value = "'0.07168"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("D11").Value = "'0.8861"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.28%  "
$ws.Range("D12").Value = "1.885.93"
$ws.Range("E12").Value = "  +1.75%  "
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").Value = "'5.348"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.14%  "
$ws.Range("D15").Value = "'89.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.88%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "'0.000008576"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("D19").Value = "'1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "27.290.67"
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("D21").Value = "'5.045"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.02%  "
$ws.Range("D22").Value = "2.136.08"
$ws.Range("E22").Value = "  +2.01%  "
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("D24").Value = "'6.487"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'151.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("D26").Value = "'1.848"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").Value = "'18.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("D28").Value = "'2.172"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.20%  "
$ws.Range("D29").Value = "'112.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("D30").Value = "'4.759"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.86%  "
$ws.Range("D31").Value = "'4.712"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("D32").Value = "'0.09054"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("D33").Value = "'0.05191"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.00%  "
$ws.Range("D34").Value = "'3.116"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.95%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7571"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.177"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.21%  "
$ws.Range("D37").Value = "'0.02046"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "'2.541"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").Value = "'3.046"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("D40").Value = "'1.085"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("D41").Value = "'0.5456"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("D42").Value = "'6.690"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.64%  "
$ws.Range("D43").Value = "'115.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.61%  "
$ws.Range("D44").Value = "'8.518"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("D45").Value = "'0.1491"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.90%  "
$ws.Range("D46").Value = "'0.4705"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("E47").Value = "  -3.50%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").Value = "'1.580"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.98%  "
$ws.Range("D50").Value = "'65.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.58%  "
$ws.Range("D51").Value = "'36.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.02%  "
